# [Kadastro App] Yeni kayit eklendi: 120 - 03.08.2025 22:11:19
#
# Adds a new record row for "120" dated 2025-08-03 to both the master
# "Kayitlar" log (sheet1, next free row after the existing data) and the
# matching district sheet "Merkez Ilce" (sheet5), keeping every field
# stored as plain text (the same convention used by all the other rows).

$wb = $excel.ActiveWorkbook

$cols = @("A", "B", "C", "D", "E", "F", "G")
$values = @(
    "120",
    "2025-08-03",
    "Merkez İlçe",
    "3",
    "2",
    "DÜZELTME",
    "HİKMET GÜLCAN (K.Mühendisi)"
)

# sheet name -> row index to write the new record into
$targets = @{
    "Kayitlar"    = 19
    "Merkez İlçe" = 2
}

foreach ($sheetName in $targets.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowIndex = $targets[$sheetName]

    for ($i = 0; $i -lt $cols.Length; $i++) {
        $cell = $ws.Range($cols[$i] + $rowIndex)
        # Force text storage (like the rest of the sheet) so numeric-looking
        # values such as "120", "3" or "2" are not reinterpreted as numbers,
        # and date-looking values such as "2025-08-03" are not reinterpreted
        # as dates.
        $cell.NumberFormat = "@"
        $cell.Value = $values[$i]
        # Drop back to the default "Normal" style once the value has been
        # committed as text so we don't leave a stray number format behind.
        $cell.Style = "Normal"
    }
}
